# Refresh cryptocurrency price/volume snapshot values (and a handful of
# rows whose ranking order swapped) to match the latest GitHub Actions run.
#
# Numeric-looking strings (e.g. "556.95") are written with a leading
# apostrophe to force Excel to keep them as text instead of silently
# converting them to floating point numbers (which would both lose the
# original text formatting, e.g. "0.130" -> 0.13, and introduce binary
# rounding noise, e.g. "556.95" -> 556.95000000000005). The Style is then
# reset to "Normal" so the quote-prefix text hint doesn't leave a stray
# cell style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.469.55"
$ws.Range("E2").Value = "  -5.41%  "
$ws.Range("D3").Value = "3.353.48"
$ws.Range("E3").Value = "  -6.31%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'556.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.73%  "
$ws.Range("D6").Value = "'182.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.53%  "
$ws.Range("E7").Value = "  -4.76%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "3.344.41"
$ws.Range("E9").Value = "  -6.21%  "
$ws.Range("E10").Value = "  -12.87%  "
$ws.Range("D11").Value = "'0.591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.88%  "
$ws.Range("D12").Value = "'47.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.41%  "
$ws.Range("E13").Value = "  -10.69%  "
$ws.Range("E14").Value = "  -9.52%  "
$ws.Range("D15").Value = "3.883.33"
$ws.Range("D16").Value = "'597.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -13.83%  "
$ws.Range("D17").Value = "66.336.36"
$ws.Range("E17").Value = "  -5.72%  "
$ws.Range("D18").Value = "3.350.40"
$ws.Range("E18").Value = "  -6.73%  "
$ws.Range("E19").Value = "  -4.56%  "
$ws.Range("D20").Value = "'17.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.46%  "
$ws.Range("D21").Value = "'11.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.44%  "
$ws.Range("E22").Value = "  -8.35%  "
$ws.Range("D23").Value = "'16.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.40%  "
$ws.Range("E24").Value = "  -5.22%  "
$ws.Range("D25").Value = "'96.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -12.43%  "
$ws.Range("D26").Value = "'4.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.78%  "
$ws.Range("D27").Value = "'2.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.50%  "
$ws.Range("D28").Value = "'9.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.49%  "
$ws.Range("D29").Value = "'8.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.90%  "
$ws.Range("D30").Value = "'30.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.14%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.47%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "'3.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -13.22%  "
$ws.Range("D33").Value = "'11.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.60%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.105"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.60%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.821.59"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "'533.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.64%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'57.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.70%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'3.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +33.46%  "
$ws.Range("D40").Value = "'3.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.64%  "
$ws.Range("D41").Value = "0.0₃0721"
$ws.Range("E41").Value = "  -14.23%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.99%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.128"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.97%  "
$ws.Range("E44").Value = "  -8.56%  "
$ws.Range("D45").Value = "'32.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.88%  "
$ws.Range("D47").Value = "'2.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.01%  "
$ws.Range("D48").Value = "'3.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.70%  "
$ws.Range("D49").Value = "'0.130"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.44%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("E51").Value = "  -10.79%  "
